$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.87"
$ws.Range("E2").Value = "'-3.69%"
$ws.Range("D3").Value = "'37.02"
$ws.Range("E3").Value = "'-6.97%"
$ws.Range("D4").Value = "'5.093"
$ws.Range("E4").Value = "'-0.95%"
$ws.Range("D5").Value = "'0.07720"
$ws.Range("E5").Value = "'-6.11%"
$ws.Range("D6").Value = "'4.361"
$ws.Range("E6").Value = "'0.66%"
$ws.Range("D7").Value = "'8.207"
$ws.Range("E7").Value = "'-1.83%"
$ws.Range("D8").Value = "'1.876"
$ws.Range("E8").Value = "'-8.97%"
$ws.Range("D10").Value = "'0.9181"
$ws.Range("E10").Value = "'-2.36%"
$ws.Range("D11").Value = "'0.1208"
$ws.Range("E11").Value = "'-11.63%"
$ws.Range("D12").Value = "'0.1888"
$ws.Range("E12").Value = "'-5.03%"
$ws.Range("D13").Value = "'0.08710"
$ws.Range("E13").Value = "'-4.89%"
$ws.Range("D14").Value = "'0.03384"
$ws.Range("E14").Value = "'-3.11%"
$ws.Range("D15").Value = "'0.09698"
$ws.Range("E15").Value = "'-1.12%"
$ws.Range("D16").Value = "'0.001365"
$ws.Range("E16").Value = "'-3.41%"
$ws.Range("D17").Value = "'0.006120"
$ws.Range("E17").Value = "'0.29%"
$ws.Range("D18").Value = "'3.563"
$ws.Range("E18").Value = "'-3.48%"
$ws.Range("E19").Value = "'-2.99%"
$ws.Range("D20").Value = "'0.1276"
$ws.Range("E20").Value = "'-2.63%"
$ws.Range("D21").Value = "'5.023"
$ws.Range("E21").Value = "'1.49%"
$ws.Range("D22").Value = "'0.2496"
$ws.Range("E22").Value = "'1.87%"
$ws.Range("D23").Value = "'0.02111"
$ws.Range("E23").Value = "'5,182.90%"
$ws.Range("D24").Value = "'0.04326"
$ws.Range("E24").Value = "'-0.68%"
$ws.Range("D25").Value = "'0.001217"
$ws.Range("E25").Value = "'-1.03%"
$ws.Range("D26").Value = "'0.004469"
$ws.Range("E26").Value = "'-7.44%"
$ws.Range("D27").Value = "'0.0001355"
$ws.Range("E27").Value = "'4.17%"
$ws.Range("D39").Value = "'0.02221"
$ws.Range("E39").Value = "'-0.66%"
$ws.Range("D40").Value = "'0.04918"
$ws.Range("E40").Value = "'-5.75%"
$ws.Range("D41").Value = "'0.007604"
$ws.Range("E41").Value = "'-2.10%"
$ws.Range("D42").Value = "'0.009830"
$ws.Range("E42").Value = "'1.50%"
$ws.Range("D43").Value = "'0.1329"
$ws.Range("E43").Value = "'-5.51%"
$ws.Range("D44").Value = "'0.001999"
$ws.Range("E44").Value = "'-2.40%"
$ws.Range("E45").Value = "'-8.87%"
$ws.Range("E46").Value = "'2.67%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.41%"
$ws.Range("D48").Value = "'0.003006"
$ws.Range("E48").Value = "'2.11%"
$ws.Range("D49").Value = "'0.001304"
$ws.Range("E49").Value = "'-22.76%"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.41%"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.41%"
